# Add 17 new feed rows to the "Feeds" sheet (architectural subreddit feeds
# + one Youtube feed), matching the order data was entered: row 1105 filled
# in full, then rows 1106-1121 filled column by column (A, D, C, B, E, F) -
# this mirrors how the shared strings table was populated in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feeds")

# Row 1105: Tanmay Bhat - Youtube
$ws.Cells.Item(1105, 1).Value = "Tanmay Bhat - Youtube"
$ws.Cells.Item(1105, 2).Value = "popculture"
$ws.Cells.Item(1105, 3).Value = "http://www.youtube.com/feeds/videos.xml?channel_id=UC0rE2qq81of4fojo-KhO5rg"
$ws.Cells.Item(1105, 4).Value = "https://www.youtube.com/channel/UC0rE2qq81of4fojo-KhO5rg"
$ws.Cells.Item(1105, 5).Value = "youtube"
$ws.Cells.Item(1105, 6).Value = "personal"

# Rows 1106-1121: architectural subreddit feeds
# Column A - display titles
$ws.Cells.Item(1106, 1).Value = "r/urbanism"
$ws.Cells.Item(1107, 1).Value = "r/architectureporn"
$ws.Cells.Item(1108, 1).Value = "r/architects"
$ws.Cells.Item(1109, 1).Value = "r/exteriordesign"
$ws.Cells.Item(1110, 1).Value = "r/cityporn"
$ws.Cells.Item(1111, 1).Value = "r/urbanhell"
$ws.Cells.Item(1112, 1).Value = "r/skyscrapers"
$ws.Cells.Item(1113, 1).Value = "r/centuryhomes"
$ws.Cells.Item(1114, 1).Value = "r/architecturalrevival"
$ws.Cells.Item(1115, 1).Value = "r/interiordesign"
$ws.Cells.Item(1116, 1).Value = "r/architectureportfolio"
$ws.Cells.Item(1117, 1).Value = "r/amazing_architecture"
$ws.Cells.Item(1118, 1).Value = "r/brutalism"
$ws.Cells.Item(1119, 1).Value = "r/artdeco"
$ws.Cells.Item(1120, 1).Value = "r/floorplan"
$ws.Cells.Item(1121, 1).Value = "r/urbandesign"

# Column D - website urls
$ws.Cells.Item(1106, 4).Value = "https://www.reddit.com/r/urbanism/"
$ws.Cells.Item(1107, 4).Value = "https://www.reddit.com/r/architectureporn/"
$ws.Cells.Item(1108, 4).Value = "https://www.reddit.com/r/architects/"
$ws.Cells.Item(1109, 4).Value = "https://www.reddit.com/r/exteriordesign/"
$ws.Cells.Item(1110, 4).Value = "https://www.reddit.com/r/cityporn/"
$ws.Cells.Item(1111, 4).Value = "https://www.reddit.com/r/urbanhell/"
$ws.Cells.Item(1112, 4).Value = "https://www.reddit.com/r/skyscrapers/"
$ws.Cells.Item(1113, 4).Value = "https://www.reddit.com/r/centuryhomes/"
$ws.Cells.Item(1114, 4).Value = "https://www.reddit.com/r/architecturalrevival/"
$ws.Cells.Item(1115, 4).Value = "https://www.reddit.com/r/interiordesign/"
$ws.Cells.Item(1116, 4).Value = "https://www.reddit.com/r/architectureportfolio/"
$ws.Cells.Item(1117, 4).Value = "https://www.reddit.com/r/amazing_architecture/"
$ws.Cells.Item(1118, 4).Value = "https://www.reddit.com/r/brutalism/"
$ws.Cells.Item(1119, 4).Value = "https://www.reddit.com/r/artdeco/"
$ws.Cells.Item(1120, 4).Value = "https://www.reddit.com/r/floorplan/"
$ws.Cells.Item(1121, 4).Value = "https://www.reddit.com/r/urbandesign/"

# Column C - rss feed urls
$ws.Cells.Item(1106, 3).Value = "https://www.reddit.com/r/urbanism.rss"
$ws.Cells.Item(1107, 3).Value = "https://www.reddit.com/r/architectureporn.rss"
$ws.Cells.Item(1108, 3).Value = "https://www.reddit.com/r/architects.rss"
$ws.Cells.Item(1109, 3).Value = "https://www.reddit.com/r/exteriordesign.rss"
$ws.Cells.Item(1110, 3).Value = "https://www.reddit.com/r/cityporn.rss"
$ws.Cells.Item(1111, 3).Value = "https://www.reddit.com/r/urbanhell.rss"
$ws.Cells.Item(1112, 3).Value = "https://www.reddit.com/r/skyscrapers.rss"
$ws.Cells.Item(1113, 3).Value = "https://www.reddit.com/r/centuryhomes.rss"
$ws.Cells.Item(1114, 3).Value = "https://www.reddit.com/r/architecturalrevival.rss"
$ws.Cells.Item(1115, 3).Value = "https://www.reddit.com/r/interiordesign.rss"
$ws.Cells.Item(1116, 3).Value = "https://www.reddit.com/r/architectureportfolio.rss"
$ws.Cells.Item(1117, 3).Value = "https://www.reddit.com/r/amazing_architecture.rss"
$ws.Cells.Item(1118, 3).Value = "https://www.reddit.com/r/brutalism.rss"
$ws.Cells.Item(1119, 3).Value = "https://www.reddit.com/r/artdeco.rss"
$ws.Cells.Item(1120, 3).Value = "https://www.reddit.com/r/floorplan.rss"
$ws.Cells.Item(1121, 3).Value = "https://www.reddit.com/r/urbandesign.rss"

# Column B - topic codes
$ws.Cells.Item(1106, 2).Value = "urbanplanning"
$ws.Cells.Item(1107, 2).Value = "architecture"
$ws.Cells.Item(1108, 2).Value = "architecture"
$ws.Cells.Item(1109, 2).Value = "architecture"
$ws.Cells.Item(1110, 2).Value = "urbanplanning"
$ws.Cells.Item(1111, 2).Value = "urbanplanning"
$ws.Cells.Item(1112, 2).Value = "architecture"
$ws.Cells.Item(1113, 2).Value = "archhistory"
$ws.Cells.Item(1114, 2).Value = "architecture"
$ws.Cells.Item(1115, 2).Value = "interiordesign"
$ws.Cells.Item(1116, 2).Value = "architecture"
$ws.Cells.Item(1117, 2).Value = "architecture"
$ws.Cells.Item(1118, 2).Value = "archmovements"
$ws.Cells.Item(1119, 2).Value = "archmovements"
$ws.Cells.Item(1120, 2).Value = "architecture"
$ws.Cells.Item(1121, 2).Value = "urbanplanning"

# Column E - source_type
$ws.Cells.Item(1106, 5).Value = "reddit"
$ws.Cells.Item(1107, 5).Value = "reddit"
$ws.Cells.Item(1108, 5).Value = "reddit"
$ws.Cells.Item(1109, 5).Value = "reddit"
$ws.Cells.Item(1110, 5).Value = "reddit"
$ws.Cells.Item(1111, 5).Value = "reddit"
$ws.Cells.Item(1112, 5).Value = "reddit"
$ws.Cells.Item(1113, 5).Value = "reddit"
$ws.Cells.Item(1114, 5).Value = "reddit"
$ws.Cells.Item(1115, 5).Value = "reddit"
$ws.Cells.Item(1116, 5).Value = "reddit"
$ws.Cells.Item(1117, 5).Value = "reddit"
$ws.Cells.Item(1118, 5).Value = "reddit"
$ws.Cells.Item(1119, 5).Value = "reddit"
$ws.Cells.Item(1120, 5).Value = "reddit"
$ws.Cells.Item(1121, 5).Value = "reddit"

# Column F - owner_type
$ws.Cells.Item(1106, 6).Value = "organization"
$ws.Cells.Item(1107, 6).Value = "organization"
$ws.Cells.Item(1108, 6).Value = "organization"
$ws.Cells.Item(1109, 6).Value = "organization"
$ws.Cells.Item(1110, 6).Value = "organization"
$ws.Cells.Item(1111, 6).Value = "organization"
$ws.Cells.Item(1112, 6).Value = "organization"
$ws.Cells.Item(1113, 6).Value = "organization"
$ws.Cells.Item(1114, 6).Value = "organization"
$ws.Cells.Item(1115, 6).Value = "organization"
$ws.Cells.Item(1116, 6).Value = "organization"
$ws.Cells.Item(1117, 6).Value = "organization"
$ws.Cells.Item(1118, 6).Value = "organization"
$ws.Cells.Item(1119, 6).Value = "organization"
$ws.Cells.Item(1120, 6).Value = "organization"
$ws.Cells.Item(1121, 6).Value = "organization"

# Reflect the final cursor position left by the edit session
[void]$ws.Range("A1107").Select()
